# Insert 3 new weekly rows for "Fruta / Vega Monumental Concepción - Frutilla"
# at the top of the date-ordered block (row 395), pushing the existing rows
# (previously 395-439) down to 398-442.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 395.
$ws.Range("A395:A397").EntireRow.Insert()

# Seed the three new rows by duplicating the rows that now sit right below
# them (which carry the same Mercado/Region/Producto/Calidad shape), then
# overwrite the date + volume/price columns with this week's values.
$ws.Rows.Item(398).Copy()
$ws.Rows.Item(395).PasteSpecial()

$ws.Rows.Item(399).Copy()
$ws.Rows.Item(396).PasteSpecial()

$ws.Rows.Item(400).Copy()
$ws.Rows.Item(397).PasteSpecial()

$excel.CutCopyMode = 0

# Row 395 - Calidad "Especial"
$ws.Range("D395").Value = 44918
$ws.Range("M395").Value = 200
$ws.Range("N395").Value = 8000
$ws.Range("O395").Value = 8000
$ws.Range("P395").Value = 8000
$ws.Range("S395").Value = 1143

# Row 396 - Calidad "Primera"
$ws.Range("D396").Value = 44918
$ws.Range("M396").Value = 100
$ws.Range("N396").Value = 7000
$ws.Range("O396").Value = 7000
$ws.Range("P396").Value = 7000
$ws.Range("S396").Value = 1000

# Row 397 - Calidad "Segunda"
$ws.Range("D397").Value = 44918
$ws.Range("M397").Value = 50
$ws.Range("N397").Value = 5000
$ws.Range("O397").Value = 5000
$ws.Range("P397").Value = 5000
$ws.Range("S397").Value = 714
